$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update goods price in seed data
$ws.Range("E4").Value = 0.02
$ws.Range("E5").Value = 0.01

# Update active selection to match the authored edit
$ws.Range("E4").Select()
